$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'25.872.08"
$ws.Cells.Item(2, 5).Value = "'  -0.14%  "

$ws.Cells.Item(3, 4).Value = "'1.642.74"
$ws.Cells.Item(3, 5).Value = "'  +0.55%  "

$ws.Cells.Item(4, 4).Value = "'1.006"
$ws.Cells.Item(4, 5).Value = "'  -1.12%  "

$ws.Cells.Item(5, 4).Value = "'215.92"
$ws.Cells.Item(5, 5).Value = "'  +0.21%  "

$ws.Cells.Item(6, 4).Value = "'0.5045"
$ws.Cells.Item(6, 5).Value = "'  -0.05%  "

$ws.Cells.Item(7, 4).Value = "'1.010"
$ws.Cells.Item(7, 5).Value = "'  -0.72%  "

$ws.Cells.Item(8, 4).Value = "'0.2571"
$ws.Cells.Item(8, 5).Value = "'  -0.18%  "

$ws.Cells.Item(9, 4).Value = "'0.06386"
$ws.Cells.Item(9, 5).Value = "'  -0.24%  "

$ws.Cells.Item(10, 4).Value = "'19.48"
$ws.Cells.Item(10, 5).Value = "'  +0.12%  "

$ws.Cells.Item(11, 4).Value = "'0.07829"
$ws.Cells.Item(11, 5).Value = "'  +0.93%  "

$ws.Cells.Item(12, 4).Value = "'1.650.70"
$ws.Cells.Item(12, 5).Value = "'  +0.72%  "

$ws.Cells.Item(13, 4).Value = "'4.265"
$ws.Cells.Item(13, 5).Value = "'  +0.25%  "

$ws.Cells.Item(14, 4).Value = "'1.862.03"
$ws.Cells.Item(14, 5).Value = "'  -0.03%  "

$ws.Cells.Item(15, 4).Value = "'0.5435"
$ws.Cells.Item(15, 5).Value = "'  -0.06%  "

$ws.Cells.Item(16, 4).Value = "'0.0₅7877"
$ws.Cells.Item(16, 5).Value = "'  -0.92%  "

$ws.Cells.Item(17, 4).Value = "'64.59"
$ws.Cells.Item(17, 5).Value = "'  +1.75%  "

$ws.Cells.Item(18, 4).Value = "'25.935.10"
$ws.Cells.Item(18, 5).Value = "'  -0.02%  "

$ws.Cells.Item(19, 4).Value = "'1.011"
$ws.Cells.Item(19, 5).Value = "'  -0.63%  "

$ws.Cells.Item(20, 4).Value = "'195.20"
$ws.Cells.Item(20, 5).Value = "'  -4.19%  "

$ws.Cells.Item(21, 4).Value = "'4.376"
$ws.Cells.Item(21, 5).Value = "'  +1.65%  "

$ws.Cells.Item(22, 4).Value = "'9.911"
$ws.Cells.Item(22, 5).Value = "'  -0.74%  "

$ws.Cells.Item(23, 4).Value = "'5.968"
$ws.Cells.Item(23, 5).Value = "'  +0.12%  "

$ws.Cells.Item(24, 4).Value = "'1.009"
$ws.Cells.Item(24, 5).Value = "'  -0.93%  "

$ws.Cells.Item(25, 4).Value = "'1.910"
$ws.Cells.Item(25, 5).Value = "'  -3.18%  "

$ws.Cells.Item(26, 4).Value = "'140.73"
$ws.Cells.Item(26, 5).Value = "'  -0.85%  "

$ws.Cells.Item(27, 4).Value = "'0.1131"
$ws.Cells.Item(27, 5).Value = "'  -1.93%  "

$ws.Cells.Item(28, 4).Value = "'6.805"
$ws.Cells.Item(28, 5).Value = "'  +0.04%  "

$ws.Cells.Item(29, 4).Value = "'15.62"
$ws.Cells.Item(29, 5).Value = "'  -0.62%  "

$ws.Cells.Item(30, 4).Value = "'1.247"
$ws.Cells.Item(30, 5).Value = "'  +0.44%  "

$ws.Cells.Item(31, 4).Value = "'0.04865"
$ws.Cells.Item(31, 5).Value = "'  -2.29%  "

$ws.Cells.Item(32, 4).Value = "'3.248"
$ws.Cells.Item(32, 5).Value = "'  -0.32%  "

$ws.Cells.Item(33, 4).Value = "'3.182"
$ws.Cells.Item(33, 5).Value = "'  -0.35%  "

$ws.Cells.Item(34, 4).Value = "'1.535"
$ws.Cells.Item(34, 5).Value = "'  -0.03%  "

$ws.Cells.Item(35, 4).Value = "'2.388"
$ws.Cells.Item(35, 5).Value = "'  +1.54%  "

$ws.Cells.Item(36, 4).Value = "'0.8891"
$ws.Cells.Item(36, 5).Value = "'  -0.08%  "

$ws.Cells.Item(37, 4).Value = "'2.611"
$ws.Cells.Item(37, 5).Value = "'  -0.62%  "

$ws.Cells.Item(38, 4).Value = "'1.131.18"
$ws.Cells.Item(38, 5).Value = "'  +1.34%  "

$ws.Cells.Item(39, 4).Value = "'0.5507"
$ws.Cells.Item(39, 5).Value = "'  -2.46%  "

$ws.Cells.Item(40, 4).Value = "'0.01562"
$ws.Cells.Item(40, 5).Value = "'  -0.25%  "

$ws.Cells.Item(41, 4).Value = "'1.013"
$ws.Cells.Item(41, 5).Value = "'  -0.42%  "

$ws.Cells.Item(42, 4).Value = "'5.683"
$ws.Cells.Item(42, 5).Value = "'  +1.08%  "

$ws.Cells.Item(43, 4).Value = "'0.8133"
$ws.Cells.Item(43, 5).Value = "'  -0.31%  "

$ws.Cells.Item(44, 4).Value = "'99.67"
$ws.Cells.Item(44, 5).Value = "'  +0.07%  "

$ws.Cells.Item(45, 2).Value = "RocketPoolETH"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(45, 4).Value = "'1.772.62"
$ws.Cells.Item(45, 5).Value = "'  -0.10%  "

$ws.Cells.Item(46, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(46, 4).Value = "'0.0₈117"
$ws.Cells.Item(46, 5).Value = "'  +2.90%  "

$ws.Cells.Item(47, 4).Value = "'0.4559"
$ws.Cells.Item(47, 5).Value = "'  +0.05%  "

$ws.Cells.Item(48, 4).Value = "'1.004"
$ws.Cells.Item(48, 5).Value = "'  -1.17%  "

$ws.Cells.Item(49, 4).Value = "'54.99"
$ws.Cells.Item(49, 5).Value = "'  +0.48%  "

$ws.Cells.Item(50, 4).Value = "'0.05046"
$ws.Cells.Item(50, 5).Value = "'  +0.18%  "

$ws.Cells.Item(51, 4).Value = "'1.010"
$ws.Cells.Item(51, 5).Value = "'  -0.63%  "
